$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,4).NumberFormat = '@'
$ws.Cells.Item(2,4).Value = '27.205.83'
$ws.Cells.Item(2,4).Style = 'Normal'
$ws.Cells.Item(2,5).Value = '  -1.90%  '

# Row 3
$ws.Cells.Item(3,4).NumberFormat = '@'
$ws.Cells.Item(3,4).Value = '1.819.02'
$ws.Cells.Item(3,4).Style = 'Normal'
$ws.Cells.Item(3,5).Value = '  -2.24%  '

# Row 4
$ws.Cells.Item(4,4).NumberFormat = '@'
$ws.Cells.Item(4,4).Value = '1.006'
$ws.Cells.Item(4,4).Style = 'Normal'
$ws.Cells.Item(4,5).Value = '  -1.62%  '

# Row 5
$ws.Cells.Item(5,4).NumberFormat = '@'
$ws.Cells.Item(5,4).Value = '314.08'
$ws.Cells.Item(5,4).Style = 'Normal'
$ws.Cells.Item(5,5).Value = '  -2.19%  '

# Row 6
$ws.Cells.Item(6,5).Value = '  -1.65%  '

# Row 7
$ws.Cells.Item(7,4).NumberFormat = '@'
$ws.Cells.Item(7,4).Value = '0.4256'
$ws.Cells.Item(7,4).Style = 'Normal'
$ws.Cells.Item(7,5).Value = '  -2.54%  '

# Row 8
$ws.Cells.Item(8,4).NumberFormat = '@'
$ws.Cells.Item(8,4).Value = '0.3665'
$ws.Cells.Item(8,4).Style = 'Normal'
$ws.Cells.Item(8,5).Value = '  -3.09%  '

# Row 9
$ws.Cells.Item(9,4).NumberFormat = '@'
$ws.Cells.Item(9,4).Value = '45.82'
$ws.Cells.Item(9,4).Style = 'Normal'
$ws.Cells.Item(9,5).Value = '  -2.42%  '

# Row 10
$ws.Cells.Item(10,4).NumberFormat = '@'
$ws.Cells.Item(10,4).Value = '0.07205'
$ws.Cells.Item(10,4).Style = 'Normal'
$ws.Cells.Item(10,5).Value = '  -2.83%  '

# Row 11
$ws.Cells.Item(11,4).NumberFormat = '@'
$ws.Cells.Item(11,4).Value = '0.8589'
$ws.Cells.Item(11,4).Style = 'Normal'
$ws.Cells.Item(11,5).Value = '  -2.70%  '

# Row 12
$ws.Cells.Item(12,4).NumberFormat = '@'
$ws.Cells.Item(12,4).Value = '20.92'
$ws.Cells.Item(12,4).Style = 'Normal'
$ws.Cells.Item(12,5).Value = '  -3.16%  '

# Row 13
$ws.Cells.Item(13,4).NumberFormat = '@'
$ws.Cells.Item(13,4).Value = '1.807.19'
$ws.Cells.Item(13,4).Style = 'Normal'
$ws.Cells.Item(13,5).Value = '  -2.26%  '

# Row 14
$ws.Cells.Item(14,4).NumberFormat = '@'
$ws.Cells.Item(14,4).Value = '6.652'
$ws.Cells.Item(14,4).Style = 'Normal'
$ws.Cells.Item(14,5).Value = '  -1.56%  '

# Row 15
$ws.Cells.Item(15,4).NumberFormat = '@'
$ws.Cells.Item(15,4).Value = '0.07108'
$ws.Cells.Item(15,4).Style = 'Normal'
$ws.Cells.Item(15,5).Value = '  -0.50%  '

# Row 16
$ws.Cells.Item(16,4).NumberFormat = '@'
$ws.Cells.Item(16,4).Value = '5.288'
$ws.Cells.Item(16,4).Style = 'Normal'
$ws.Cells.Item(16,5).Value = '  -3.72%  '

# Row 17
$ws.Cells.Item(17,4).NumberFormat = '@'
$ws.Cells.Item(17,4).Value = '87.86'
$ws.Cells.Item(17,4).Style = 'Normal'
$ws.Cells.Item(17,5).Value = '  +0.66%  '

# Row 18
$ws.Cells.Item(18,4).NumberFormat = '@'
$ws.Cells.Item(18,4).Value = '1.006'
$ws.Cells.Item(18,4).Style = 'Normal'
$ws.Cells.Item(18,5).Value = '  -1.91%  '

# Row 19
$ws.Cells.Item(19,4).NumberFormat = '@'
$ws.Cells.Item(19,4).Value = '0.000008839'
$ws.Cells.Item(19,4).Style = 'Normal'
$ws.Cells.Item(19,5).Value = '  -2.57%  '

# Row 20
$ws.Cells.Item(20,4).NumberFormat = '@'
$ws.Cells.Item(20,4).Value = '1.003'
$ws.Cells.Item(20,4).Style = 'Normal'
$ws.Cells.Item(20,5).Value = '  -1.59%  '

# Row 21
$ws.Cells.Item(21,4).NumberFormat = '@'
$ws.Cells.Item(21,4).Value = '15.00'
$ws.Cells.Item(21,4).Style = 'Normal'
$ws.Cells.Item(21,5).Value = '  -2.98%  '

# Row 22
$ws.Cells.Item(22,4).NumberFormat = '@'
$ws.Cells.Item(22,4).Value = '27.236.25'
$ws.Cells.Item(22,4).Style = 'Normal'
$ws.Cells.Item(22,5).Value = '  -1.81%  '

# Row 23
$ws.Cells.Item(23,4).NumberFormat = '@'
$ws.Cells.Item(23,4).Value = '5.131'
$ws.Cells.Item(23,4).Style = 'Normal'
$ws.Cells.Item(23,5).Value = '  -3.06%  '

# Row 24
$ws.Cells.Item(24,4).NumberFormat = '@'
$ws.Cells.Item(24,4).Value = '10.86'
$ws.Cells.Item(24,4).Style = 'Normal'
$ws.Cells.Item(24,5).Value = '  -2.63%  '

# Row 25
$ws.Cells.Item(25,4).NumberFormat = '@'
$ws.Cells.Item(25,4).Value = '2.048.49'
$ws.Cells.Item(25,4).Style = 'Normal'
$ws.Cells.Item(25,5).Value = '  -2.73%  '

# Row 26
$ws.Cells.Item(26,4).NumberFormat = '@'
$ws.Cells.Item(26,4).Value = '2.005'
$ws.Cells.Item(26,4).Style = 'Normal'

# Row 27
$ws.Cells.Item(27,4).NumberFormat = '@'
$ws.Cells.Item(27,4).Value = '153.07'
$ws.Cells.Item(27,4).Style = 'Normal'
$ws.Cells.Item(27,5).Value = '  -2.66%  '

# Row 28
$ws.Cells.Item(28,4).NumberFormat = '@'
$ws.Cells.Item(28,4).Value = '18.27'
$ws.Cells.Item(28,4).Style = 'Normal'
$ws.Cells.Item(28,5).Value = '  -2.50%  '

# Row 29
$ws.Cells.Item(29,4).NumberFormat = '@'
$ws.Cells.Item(29,4).Value = '2.101'
$ws.Cells.Item(29,4).Style = 'Normal'
$ws.Cells.Item(29,5).Value = '  +5.51%  '

# Row 30
$ws.Cells.Item(30,4).NumberFormat = '@'
$ws.Cells.Item(30,4).Value = '5.209'
$ws.Cells.Item(30,4).Style = 'Normal'
$ws.Cells.Item(30,5).Value = '  -3.40%  '

# Row 31
$ws.Cells.Item(31,4).NumberFormat = '@'
$ws.Cells.Item(31,4).Value = '115.95'
$ws.Cells.Item(31,4).Style = 'Normal'
$ws.Cells.Item(31,5).Value = '  -4.61%  '

# Row 32
$ws.Cells.Item(32,4).NumberFormat = '@'
$ws.Cells.Item(32,4).Value = '0.08862'
$ws.Cells.Item(32,4).Style = 'Normal'
$ws.Cells.Item(32,5).Value = '  -2.21%  '

# Row 33
$ws.Cells.Item(33,2).Value = 'ARBITRUM'
$ws.Cells.Item(33,3).Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Cells.Item(33,4).NumberFormat = '@'
$ws.Cells.Item(33,4).Value = '1.189'
$ws.Cells.Item(33,4).Style = 'Normal'
$ws.Cells.Item(33,5).Value = '  -2.17%  '

# Row 34
$ws.Cells.Item(34,2).Value = 'ImmutableX'
$ws.Cells.Item(34,3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Cells.Item(34,4).NumberFormat = '@'
$ws.Cells.Item(34,4).Value = '0.7575'
$ws.Cells.Item(34,4).Style = 'Normal'
$ws.Cells.Item(34,5).Value = '  -1.53%  '

# Row 35
$ws.Cells.Item(35,4).NumberFormat = '@'
$ws.Cells.Item(35,4).Value = '4.438'
$ws.Cells.Item(35,4).Style = 'Normal'
$ws.Cells.Item(35,5).Value = '  -2.64%  '

# Row 36
$ws.Cells.Item(36,4).NumberFormat = '@'
$ws.Cells.Item(36,4).Value = '2.803'
$ws.Cells.Item(36,4).Style = 'Normal'
$ws.Cells.Item(36,5).Value = '  -7.79%  '

# Row 37
$ws.Cells.Item(37,4).NumberFormat = '@'
$ws.Cells.Item(37,4).Value = '1.004'
$ws.Cells.Item(37,4).Style = 'Normal'
$ws.Cells.Item(37,5).Value = '  -1.74%  '

# Row 38
$ws.Cells.Item(38,5).Value = '  -2.54%  '

# Row 39
$ws.Cells.Item(39,4).NumberFormat = '@'
$ws.Cells.Item(39,4).Value = '0.01957'
$ws.Cells.Item(39,4).Style = 'Normal'
$ws.Cells.Item(39,5).Value = '  -0.94%  '

# Row 40
$ws.Cells.Item(40,4).NumberFormat = '@'
$ws.Cells.Item(40,4).Value = '0.05248'
$ws.Cells.Item(40,4).Style = 'Normal'
$ws.Cells.Item(40,5).Value = '  -1.15%  '

# Row 41
$ws.Cells.Item(41,4).NumberFormat = '@'
$ws.Cells.Item(41,4).Value = '2.897'
$ws.Cells.Item(41,4).Style = 'Normal'
$ws.Cells.Item(41,5).Value = '  +0.54%  '

# Row 42
$ws.Cells.Item(42,4).NumberFormat = '@'
$ws.Cells.Item(42,4).Value = '7.083'
$ws.Cells.Item(42,4).Style = 'Normal'
$ws.Cells.Item(42,5).Value = '  +1.48%  '

# Row 43
$ws.Cells.Item(43,4).NumberFormat = '@'
$ws.Cells.Item(43,4).Value = '0.1673'
$ws.Cells.Item(43,4).Style = 'Normal'
$ws.Cells.Item(43,5).Value = '  -0.25%  '

# Row 45
$ws.Cells.Item(45,4).NumberFormat = '@'
$ws.Cells.Item(45,4).Value = '8.582'
$ws.Cells.Item(45,4).Style = 'Normal'
$ws.Cells.Item(45,5).Value = '  -1.60%  '

# Row 46
$ws.Cells.Item(46,4).NumberFormat = '@'
$ws.Cells.Item(46,4).Value = '10.53'
$ws.Cells.Item(46,4).Style = 'Normal'
$ws.Cells.Item(46,5).Value = '  -1.70%  '

# Row 47
$ws.Cells.Item(47,4).NumberFormat = '@'
$ws.Cells.Item(47,4).Value = '106.45'
$ws.Cells.Item(47,4).Style = 'Normal'
$ws.Cells.Item(47,5).Value = '  -3.41%  '

# Row 48
$ws.Cells.Item(48,4).NumberFormat = '@'
$ws.Cells.Item(48,4).Value = '0.4687'
$ws.Cells.Item(48,4).Style = 'Normal'
$ws.Cells.Item(48,5).Value = '  -0.89%  '

# Row 49
$ws.Cells.Item(49,4).NumberFormat = '@'
$ws.Cells.Item(49,4).Value = '1.003'
$ws.Cells.Item(49,4).Style = 'Normal'
$ws.Cells.Item(49,5).Value = '  -1.81%  '

# Row 50
$ws.Cells.Item(50,4).NumberFormat = '@'
$ws.Cells.Item(50,4).Value = '0.06387'
$ws.Cells.Item(50,4).Style = 'Normal'
$ws.Cells.Item(50,5).Value = '  -1.62%  '

# Row 51
$ws.Cells.Item(51,4).NumberFormat = '@'
$ws.Cells.Item(51,4).Value = '1.655'
$ws.Cells.Item(51,4).Style = 'Normal'
$ws.Cells.Item(51,5).Value = '  -3.49%  '
